$wb = $excel.ActiveWorkbook

# --- Sheet "ExpenseRequest": requestor/contact/event renamed ---
$wsExpense = $wb.Worksheets.Item("ExpenseRequest")
$wsExpense.Range("B2").Value = "Amanda Donovan"
$wsExpense.Range("C2").Value = "Amanda Donovan"
$wsExpense.Range("D2").Value = "PFG Golf Event"

# --- Sheet "Users": same requestor rename + widened column A ---
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("A2").Value = "Amanda Donovan"
$wsUsers.Columns.Item(1).ColumnWidth = 16.43

# --- Sheet "Approver": password updated ---
$wsApprover = $wb.Worksheets.Item("Approver")
$wsApprover.Range("B2").Value = "Bingo@1234"

# --- Restore/update on-screen selections to match the saved view state ---
$wsUsers.Activate()
$wsUsers.Range("F11").Select()

$wsExpense.Activate()
$wsExpense.Range("E12").Select()
